$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.564.53"
$ws.Range("E2").Value = "  -2.70%  "

$ws.Range("D3").Value = "1.998.22"
$ws.Range("E3").Value = "  -5.08%  "

$ws.Range("E4").Value = "  +0.87%  "

$ws.Range("D5").Value = "'330.44"
$ws.Range("E5").Value = "  -4.06%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").Value = "'0.4997"

$ws.Range("D8").Value = "'0.4240"
$ws.Range("E8").Value = "  -4.31%  "

$ws.Range("D9").Value = "'54.59"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").Value = "'0.08990"
$ws.Range("E10").Value = "  -4.16%  "

$ws.Range("D11").Value = "'1.119"
$ws.Range("E11").Value = "  -4.47%  "

$ws.Range("D12").Value = "'23.28"
$ws.Range("E12").Value = "  -6.60%  "

$ws.Range("D13").Value = "2.013.40"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D14").Value = "'8.059"
$ws.Range("E14").Value = "  -6.99%  "

$ws.Range("D15").Value = "'6.500"
$ws.Range("E15").Value = "  -6.20%  "

$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").Value = "'94.38"
$ws.Range("E17").Value = "  -7.31%  "

$ws.Range("D18").Value = "'0.00001111"
$ws.Range("E18").Value = "  -4.30%  "

$ws.Range("D19").Value = "'0.06670"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").Value = "'19.76"
$ws.Range("E20").Value = "  -6.91%  "

$ws.Range("D22").Value = "'5.954"
$ws.Range("E22").Value = "  -6.76%  "

$ws.Range("D23").Value = "29.614.14"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("D24").Value = "'12.01"
$ws.Range("E24").Value = "  -4.62%  "

$ws.Range("D25").Value = "'2.280"
$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").Value = "'159.52"
$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("D27").Value = "'20.70"
$ws.Range("E27").Value = "  -5.53%  "

$ws.Range("D28").Value = "'6.365"
$ws.Range("E28").Value = "  -5.35%  "

$ws.Range("D29").Value = "'2.305"
$ws.Range("E29").Value = "  -8.46%  "

$ws.Range("D30").Value = "'128.33"
$ws.Range("E30").Value = "  -3.93%  "

$ws.Range("D31").Value = "'1.054"
$ws.Range("E31").Value = "  -7.63%  "

$ws.Range("E32").Value = "  -5.66%  "

$ws.Range("D33").Value = "'1.567"
$ws.Range("E33").Value = "  -7.24%  "

$ws.Range("D34").Value = "'5.843"
$ws.Range("E34").Value = "  -6.41%  "

$ws.Range("D35").Value = "'3.807"
$ws.Range("E35").Value = "  -2.89%  "

$ws.Range("D36").Value = "'9.440"
$ws.Range("E36").Value = "  -8.42%  "

$ws.Range("D37").Value = "'0.02465"
$ws.Range("E37").Value = "  -6.50%  "

$ws.Range("D38").Value = "'1.310"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").Value = "'0.06350"
$ws.Range("E39").Value = "  -6.40%  "

$ws.Range("D40").Value = "'0.6574"

$ws.Range("D41").Value = "'11.66"
$ws.Range("E41").Value = "  -7.17%  "

$ws.Range("D42").Value = "'0.2053"
$ws.Range("E42").Value = "  -7.74%  "

$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("D44").Value = "'0.6333"
$ws.Range("E44").Value = "  -7.62%  "

$ws.Range("D45").Value = "'13.54"
$ws.Range("E45").Value = "  -6.83%  "

$ws.Range("D47").Value = "'1.305"
$ws.Range("E47").Value = "  -6.77%  "

$ws.Range("D48").Value = "'3.525"
$ws.Range("E48").Value = "  -3.27%  "

$ws.Range("D49").Value = "'0.00000000341"
$ws.Range("E49").Value = "  -4.71%  "

$ws.Range("D50").Value = "'0.06984"
$ws.Range("E50").Value = "  -3.62%  "

$ws.Range("D51").Value = "'1.121"
$ws.Range("E51").Value = "  -7.74%  "
